$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("neg_reaction10")

$values = @(0,2,3,4,5,8,9,11,13,15,16,18,20,22,23,24,25,27,28,29,31,32,33,34,35,36,37,39,40,41,42,44,45,46,47,48,49,50,51,52,54,55,56,57,59,60,61,62,66,67,69,70,71,72,73,74,75)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $values[$i]
}
